$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Header row (row 1): shrink the italic header-label font from 11pt (sz 22) to 10pt (sz 20)
for ($c = 1; $c -le $t.Columns.Count; $c++) {
    $cell = $t.Cell(1, $c)
    $cellStart = $cell.Range.Start
    $cellEnd = $cell.Range.End
    if (($cellEnd - 1) -gt $cellStart) {
        $r = $d.Range($cellStart, $cellEnd - 1)
        $r.Font.Size = 10
    } else {
        $r = $cell.Range.Duplicate()
        [void]$r.MoveEnd(1, -1)
        $r.Font.Size = 10
    }
}

# First column, data rows (rows 2 and 3: "Cod" and "Hake"): right-align the label paragraph
for ($r = 2; $r -le $t.Rows.Count; $r++) {
    $cell = $t.Cell($r, 1)
    $cell.Range.ParagraphFormat.Alignment = 2
}
